# Update countries & provincias Spain
# Applies the daily data refresh to the "Pais" sheet:
#  - Reorders "Haiti" so it appears before "Guyana" in the country list
#    (rows 158-160: Haiti, Guyana, San Martin (Parte Holandesa))
#  - Refreshes the statistic columns (Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
#    for the rows whose figures moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder Haiti / Guyana / San Martin (Parte Holandesa) ---
$ws.Range("A158").Value = "Haiti"
$ws.Range("A159").Value = "Guyana"
$ws.Range("A160").Value = "San Martin (Parte Holandesa)"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 986045
$ws.Range("C4").Value = 25394
$ws.Range("D4").Value = 118777
$ws.Range("E4").Value = 811891
$ws.Range("F4").Value = 15143
$ws.Range("G4").Value = 1121
$ws.Range("H4").Value = 55377

# --- Row 14: Brasil ---
$ws.Range("B14").Value = 62787
$ws.Range("C14").Value = 3591
$ws.Range("E14").Value = 28367
$ws.Range("G14").Value = 223
$ws.Range("H14").Value = 4268

# --- Row 18: Suiza ---
$ws.Range("D18").Value = 21800
$ws.Range("E18").Value = 5651

# --- Row 46: Australia ---
$ws.Range("B46").Value = 6716
$ws.Range("C46").Value = 21
$ws.Range("D46").Value = 5560
$ws.Range("E46").Value = 1073

# --- Row 51: Colombia ---
$ws.Range("B51").Value = 5379
$ws.Range("C51").Value = 237
$ws.Range("D51").Value = 1133
$ws.Range("E51").Value = 4002
$ws.Range("G51").Value = 11
$ws.Range("H51").Value = 244

# --- Row 144: Trinidad yTobago ---
$ws.Range("D144").Value = 54
$ws.Range("E144").Value = 53

# --- Row 155: Barbados ---
$ws.Range("D155").Value = 39
$ws.Range("E155").Value = 34

# --- Row 158: Haiti (new data, shares row with reordered name) ---
$ws.Range("C158").Value = 2
$ws.Range("D158").Value = 7
$ws.Range("E158").Value = 61
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 6

# --- Row 159: Guyana ---
$ws.Range("D159").Value = 12
$ws.Range("E159").Value = 54
$ws.Range("F159").Value = 5
$ws.Range("H159").Value = 8

# --- Row 160: San Martin (Parte Holandesa) ---
$ws.Range("B160").Value = 74
$ws.Range("C160").Value = 1
$ws.Range("D160").Value = 33
$ws.Range("E160").Value = 28
$ws.Range("F160").Value = 1
$ws.Range("G160").Value = 1
$ws.Range("H160").Value = 13
